$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells I1 and J1, matching the style of the existing headers (H1)
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the I2:I58 and J2:J58 data columns
$iValues = @(8,9,9,9,9,9,9,9,8,9,9,9,9,8,9,9,10,8,8,8,8,6,11,7,8,8,8,8,8,10,8,7,8,9,8,8,8,8,6,9,8,7,7,8,8,5,7,8,8,8,7,7,7,9,4,8,6)
$jValues = @(8,9,9,9,9,9,9,9,8,9,9,9,9,8,9,9,10,8,8,8,8,7,11,7,8,8,8,8,8,10,8,8,8,9,8,8,8,8,6,9,8,7,7,8,8,6,7,8,9,8,8,8,7,9,4,8,6)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}

Write-Host "Applied I0/IF columns"
